$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.549.74"
$ws.Range("E2").Value = "  -6.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.441.44"
$ws.Range("E3").Value = "  -9.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "468.07"
$ws.Range("E5").Value = "  -6.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.87"
$ws.Range("E6").Value = "  -6.00%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.493"
$ws.Range("E8").Value = "  -6.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.439.44"
$ws.Range("E9").Value = "  -9.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0951"
$ws.Range("E10").Value = "  -8.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  -12.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.315"
$ws.Range("E12").Value = "  -9.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.122"
$ws.Range("E13").Value = "  -3.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.869.42"
$ws.Range("E14").Value = "  -9.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "54.556.49"
$ws.Range("E15").Value = "  -6.74%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.64"
$ws.Range("E17").Value = "  -8.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.443.14"
$ws.Range("E18").Value = "  -9.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.22"
$ws.Range("E19").Value = "  -10.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.53"
$ws.Range("E20").Value = "  -6.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.54"
$ws.Range("E21").Value = "  -12.63%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  -13.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "56.34"
$ws.Range("E25").Value = "  -10.31%  "
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.385"
$ws.Range("E27").Value = "  -8.93%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.557.07"
$ws.Range("E28").Value = "  -9.26%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.156"
$ws.Range("E29").Value = "  -8.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.09"
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").Value = "  -12.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "145.98"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.77"
$ws.Range("E34").Value = "  -6.70%  "
$ws.Range("E35").Value = "  -9.80%  "
$ws.Range("E36").Value = "  -6.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.56"
$ws.Range("E37").Value = "  -13.84%  "
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.790"
$ws.Range("E39").Value = "  -14.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "32.91"
$ws.Range("E41").Value = "  -6.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.598"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.25"
$ws.Range("E43").Value = "  -8.12%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0522"
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.09"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("E46").Value = "  -9.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.934.41"
$ws.Range("E47").Value = "  -11.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0883"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "232.93"
$ws.Range("E50").Value = "  +6.66%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.56"
$ws.Range("E51").Value = "  -11.24%  "
